# "changed BAC to end PV"
#
# The BAC (Budget at Completion) cell on the "Kennzahlen" sheet used to hold
# a hard-coded value (1.25 * 10^6 = 1,250,000). It now references the
# cumulative PV total ("end PV") computed on the "Budgetierte Kosten" sheet,
# cell P11 (= 1,500,000). Every table / chart that derives from BAC
# (Tabelle467812, Tabelle46781213, chart1, ...) recalculates automatically.

$wb = $excel.ActiveWorkbook

$wsKennzahlen = $wb.Worksheets.Item("Kennzahlen")
$bacCell = $wsKennzahlen.Range("E92")
$bacCell.Formula = "='Budgetierte Kosten'!`$P`$11"

# Setting a formula that references a single formatted cell makes the engine
# (like Excel itself) inherit that cell's number format when the target cell
# had none. The original BAC cell is unformatted (General / no explicit
# style), so restore that by pasting formats-only from a neighboring,
# unformatted cell in the same row.
$wsKennzahlen.Range("D92").Copy() | Out-Null
$bacCell.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# Recalculate so every dependent formula (tables, charts, ...) picks up the
# new BAC-derived values.
$excel.CalculateFull()

# --- View-state bookkeeping -------------------------------------------------
# The workbook used to keep "Diagramme" (3rd tab) as the active/selected
# sheet; it now ends up on "Kennzahlen" (1st tab) with a different scroll
# position and cell selection, while "Diagramme" loses its tabSelected flag.

$wsDiagramme = $wb.Worksheets.Item("Diagramme")
$wsDiagramme.Activate()
$wsDiagramme.Range("R23").Select() | Out-Null

$wsBudget = $wb.Worksheets.Item("Budgetierte Kosten")
$wsBudget.Activate()
$wsBudget.Application.ActiveWindow.ScrollColumn = 4
$wsBudget.Range("C15").Select() | Out-Null

$wsKennzahlen.Activate()
$wsKennzahlen.Application.ActiveWindow.ScrollRow = 76
$wsKennzahlen.Application.ActiveWindow.ScrollColumn = 1
$wsKennzahlen.Range("J90").Select() | Out-Null
